$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

$ws.Range("A1").Value = "Type"

$ws.Application.Goto($ws.Range("C7"))
